# Update "want to go" counts (F column) for the 合肥·第九届环形宇宙动漫游戏嘉年华
# and 合肥·九号幻想动漫游戏嘉年华 entries, on both the "展览" sheet and the
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): row 5 -> F5, row 7 -> F7
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 5183
$wsExhibition.Range("F7").Value = 64

# Sheet "全部类型" (All types): row 9 -> F9, row 11 -> F11
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 5183
$wsAll.Range("F11").Value = 64
